$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B2 was empty, now gets the region name (text)
$ws.Range("B2").Value = "REGIONE_LAZIO"

# D2 workflowInstanceId suffix changed (ff6ce61428 -> f73afc3a33)
$ws.Range("D2").Value = "2.16.840.1.113883.2.9.2.120.4.4.b0f3ffcf25ce2aafc7dc901e2febc51f43837f4ca0fe3b6d1b02194e9047b6db.f73afc3a33^^^^urn:ihe:iti:xdw:2013:workflowInstanceId"

# E2 UAT_GTW_ID changed (1721652226262 -> 1721654256960)
$ws.Range("E2").Value = "2.16.840.1.113883.2.9.2.110.4.4^UAT_GTW_ID1721654256960"

# F2 timestamp changed (14:43:47 -> 15:17:38), keep as text
$ws.Range("F2").Value = "22-07-2024:15:17:38"
